$wb = $excel.ActiveWorkbook

# --- "About" sheet: bump the last-updated date (C1) from 3/15/2024 to 3/25/2024 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45376

# --- "RAF-capacity" sheet: hydrogen combustion turbine / combined cycle RAF raised to 1 ---
$wsCapacity = $wb.Worksheets.Item("RAF-capacity")
$wsCapacity.Range("B24").Value = 1
$wsCapacity.Range("B25").Value = 1

# --- Window/view state: RAF-capacity becomes the active tab, zoomed to 80%, with E8 selected ---
$wsCapacity.Activate()
$wsCapacity.Range("E8").Select()
$excel.ActiveWindow.Zoom = 80
